$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 tracks progress per column/topic. H4 ("Label Functionality") is now
# done, matching the rest of the row, and I4 ("Compose Mail") is in progress.
$ws.Range("H4").Value = "done"
$ws.Range("I4").Value = "in progress"

# Reflect the newly filled-in columns in the window's scroll position and
# active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 6
$ws.Range("I4").Select()
